$wb = $excel.ActiveWorkbook

# --- Sheet 1: "studyData" ---
# Add a new row 3 that duplicates row 2, with updated studyId, Pmid and ascertainment.
$ws1 = $wb.Worksheets.Item("studyData")

$ws1.Range("A3").Value = 3
$ws1.Range("B3").Value = "4000000"
$ws1.Range("C3").Value = $ws1.Range("C2").Value
$ws1.Range("D3").Value = $ws1.Range("D2").Value
$ws1.Range("E3").Value = "African"
$ws1.Range("F3").Value = $ws1.Range("F2").Value
$ws1.Range("G3").Value = $ws1.Range("G2").Value
$ws1.Range("H3").Value = $ws1.Range("H2").Value
$ws1.Range("I3").Value = 18
$ws1.Range("J3").Value = 1
$ws1.Range("K3").Value = $ws1.Range("K2").Value
$ws1.Range("L3").Value = 195
$ws1.Range("M3").Value = 36
$ws1.Range("N3").Value = 10.99
$ws1.Range("O3").Value = 196
$ws1.Range("P3").Value = 38
$ws1.Range("Q3").Value = 10.99
$ws1.Range("R3").Value = 1
$ws1.Range("S3").Value = 1
$ws1.Range("T3").Value = 1
$ws1.Range("U3").Value = 1
$ws1.Range("V3").Value = 1
$ws1.Range("W3").Value = $ws1.Range("W2").Value
$ws1.Range("X3").Value = $ws1.Range("X2").Value
$ws1.Range("Y3").Value = $ws1.Range("Y2").Value
$ws1.Range("Z3").Value = $ws1.Range("Z2").Value

# --- Sheet 2: "SNP_entryData" ---
# Update row 2 with new values and add a new row 3 mostly duplicating row 2.
$ws2 = $wb.Worksheets.Item("SNP_entryData")

$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = "xas"
$ws2.Range("C2").Value = 3.4
$ws2.Range("D2").Value = "cx1"
$ws2.Range("E2").Value = "X"
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = "1"
$ws2.Range("H2").Value = 5
$ws2.Range("I2").Value = 5
$ws2.Range("J2").Value = 5
$ws2.Range("K2").Value = "white"
$ws2.Range("L2").Value = "A"
$ws2.Range("M2").Value = "C"
$ws2.Range("N2").Value = 23
$ws2.Range("O2").Value = 23
$ws2.Range("P2").Value = 23
$ws2.Range("Q2").Value = "high"
$ws2.Range("R2").Value = "high"
$ws2.Range("S2").Value = "high.com"
$ws2.Range("T2").Value = "KROK"
$ws2.Range("U2").Value = 3
$ws2.Range("V2").Value = "high"
$ws2.Range("W2").Value = "high"
$ws2.Range("X2").Value = 1
$ws2.Range("Y2").Value = 10
$ws2.Range("Z2").Value = "high"
$ws2.Range("AA2").Value = "high.com"
$ws2.Range("AB2").Value = "reading"

$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = "mee"
$ws2.Range("C3").Value = 3.4
$ws2.Range("D3").Value = "cx1"
$ws2.Range("E3").Value = "X"
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = "1"
$ws2.Range("H3").Value = 5
$ws2.Range("I3").Value = 5
$ws2.Range("J3").Value = 5
$ws2.Range("K3").Value = "white"
$ws2.Range("L3").Value = "A"
$ws2.Range("M3").Value = "C"
$ws2.Range("N3").Value = 23
$ws2.Range("O3").Value = 23
$ws2.Range("P3").Value = 23
$ws2.Range("Q3").Value = "high"
$ws2.Range("R3").Value = "high"
$ws2.Range("S3").Value = "high.com"
$ws2.Range("T3").Value = "KROK"
$ws2.Range("U3").Value = 3
$ws2.Range("V3").Value = "high"
$ws2.Range("W3").Value = "high"
$ws2.Range("X3").Value = 1
$ws2.Range("Y3").Value = 10
$ws2.Range("Z3").Value = "high"
$ws2.Range("AA3").Value = "high.com"
$ws2.Range("AB3").Value = "learning"
